$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reviewer/name value in L2 (was "Budi", now "azara")
$ws.Range("L2").Value = "azara"

# Update the date-ish numeric value in F2
$ws.Range("F2").Value = 49998

# Update the view: select L2 as the active cell, then scroll so column E
# is the left-most visible column (topLeftCell = E1)
$ws.Range("L2").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
